$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Dyson Daniels", "PG,SG,SF", "Atlanta Hawks"),
    @("Kelly Oubre Jr.", "SG,SF", "Philadelphia 76ers"),
    @("Malik Beasley", "SG,SF", "Detroit Pistons"),
    @("Devin Vassell", "SG,SF", "San Antonio Spurs"),
    @("Jaden McDaniels", "SF,PF", "Minnesota Timberwolves"),
    @("Toumani Camara", "SG,SF,PF", "Minnesota Timberwolves"),
    @("De'Andre Hunter", "SF,PF", "Portland Trail Blazers"),
    @("Michael Porter Jr.", "SF,PF", "Cleveland Cavaliers"),
    @("Julius Randle", "PF,C", "Minnesota Timberwolves"),
    @("Guerschon Yabusele", "PF,C", "Philadelphia 76ers"),
    @("Josh Hart", "SG,SF,PF", "New York Knicks"),
    @("Donovan Mitchell", "PG,SG", "Cleveland Cavaliers"),
    @("Cam Thomas", "SG,SF", "Brooklyn Nets"),
    @("Bam Adebayo", "PF,C", "Miami Heat"),
    @("Kel'el Ware", "PF,C", "Miami Heat"),
    @("Kristaps Porzingis", "PF,C", "Boston Celtics"),
    @("Domantas Sabonis", "C", "Sacramento Kings")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Remove the now-unused last row (previously row 19) so the table shrinks
# from 18 data rows down to 17 data rows.
$ws.Range("A19:C19").Delete()
